# Apply the "complete monthly and re-run daily" update:
# For team member ID sp99069 (羅丹竺, row 17/18 on team_df; row 2 on the
# aggregated sheets), one additional arrival/count is recorded, which
# bumps the arrive_thres/count totals (and their derived ratios) up by
# one across the team_df, team_df_day, productivity_tl and
# productivity_team_function sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "team_df": rows 17, 18, 23 -> columns S (arrive_thres),
#     T (count), U (prod_hour_ratio = S / T)
$wsTeamDf = $wb.Worksheets.Item("team_df")

$teamDfRows = @(17, 18, 23)
foreach ($r in $teamDfRows) {
    $sVal = $wsTeamDf.Range("S$r").Value2
    $tVal = $wsTeamDf.Range("T$r").Value2
    $newS = $sVal + 1
    $newT = $tVal + 1
    $wsTeamDf.Range("S$r").Value = $newS
    $wsTeamDf.Range("T$r").Value = $newT
    $wsTeamDf.Range("U$r").Value = $newS / $newT
}

# --- Sheet "team_df_day": rows 2, 5, 10 -> columns F (count), G (hour-
#     equivalent denominator), H (prod_day_ratio = F / G)
$wsTeamDfDay = $wb.Worksheets.Item("team_df_day")

$teamDfDayRows = @(2, 5, 10)
foreach ($r in $teamDfDayRows) {
    $fVal = $wsTeamDfDay.Range("F$r").Value2
    $gVal = $wsTeamDfDay.Range("G$r").Value2
    $newF = $fVal + 1
    $newG = $gVal + 1
    $wsTeamDfDay.Range("F$r").Value = $newF
    $wsTeamDfDay.Range("G$r").Value = $newG
    $wsTeamDfDay.Range("H$r").Value = $newF / $newG
}

# --- Sheet "productivity_tl": rows 2, 5, 10 -> column D (TL_produtivity_score)
$wsProdTl = $wb.Worksheets.Item("productivity_tl")
$wsProdTl.Range("D2").Value = $wsTeamDfDay.Range("H2").Value2
$wsProdTl.Range("D5").Value = $wsTeamDfDay.Range("H5").Value2
$wsProdTl.Range("D10").Value = $wsTeamDfDay.Range("H10").Value2

# --- Sheet "productivity_team_function": rows 2, 5, 10 -> column D
$wsProdTeamFunc = $wb.Worksheets.Item("productivity_team_function")
$wsProdTeamFunc.Range("D2").Value = $wsTeamDfDay.Range("H2").Value2
$wsProdTeamFunc.Range("D5").Value = $wsTeamDfDay.Range("H5").Value2
$wsProdTeamFunc.Range("D10").Value = $wsTeamDfDay.Range("H10").Value2
